$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "GÊNERO" / "{{genero}}" table row entirely.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("GÊNERO", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    $cell = $rng.Cells.Item(1)
    $tbl = $cell.Range.Tables.Item(1)
    $tbl.Rows.Item($cell.RowIndex).Delete()
}

# ------------------------------------------------------------------
# 2) Insert a new, empty paragraph (same Arial/Bold/25pt/white
#    highlight formatting) right before the
#    "4. CARACTERÍSTICAS DA ÁREA DE BUSCA" heading paragraph.
# ------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("4. CARACTERÍSTICAS DA ÁREA DE BUSCA", $true, $false, $false, `
                              $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $para = $rng2.Paragraphs.Item(1)
    $para.Range.InsertParagraphBefore()
}
